# Add two new "word-frequency" sheets (聊斋志异 and 红楼梦) before 三国演义,
# each populated with A:word / B:frequency data and a treemap-style chart,
# matching the pattern already used by the other sheets in the workbook.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------
# 1. Insert the two new worksheets in the right spot / order.
#    Inserting "红楼梦" right before "三国演义" first, then inserting
#    "聊斋志异" right before "红楼梦", yields final order:
#       诗经, 唐诗三百首, 宋词三百首, 聊斋志异, 红楼梦, 三国演义, 四世同堂, 白鹿原
# ---------------------------------------------------------------------
$sanguo = $sheets.Item("三国演义")

$hlm = $sheets.Add($sanguo)
$hlm.Name = "红楼梦"

$lzzy = $sheets.Add($hlm)
$lzzy.Name = "聊斋志异"

# ---------------------------------------------------------------------
# 2. Fill in the word / frequency data (top-20 characters by frequency).
# ---------------------------------------------------------------------
$words_lzzy = @("之", "不", "曰", "人", "而", "生", "以", "一", "其", "女", "有", "为", "无", "如", "子", "见", "也", "中", "何", "者")
$freqs_lzzy = @(2.34, 1.74, 1.45, 1.34, 1.23, 1.13, 0.94, 0.9, 0.86, 0.75, 0.73, 0.63, 0.62, 0.56, 0.54, 0.54, 0.51, 0.51, 0.5, 0.5)

$words_hlm = @("了", "的", "不", "一", "来", "道", "人", "是", "说", "我", "这", "他", "你", "儿", "着", "也", "去", "玉", "有", "宝")
$freqs_hlm = @(2.93, 2.17, 2.02, 1.62, 1.56, 1.55, 1.44, 1.38, 1.33, 1.25, 1.06, 1.04, 0.98, 0.97, 0.92, 0.84, 0.84, 0.83, 0.81, 0.81)

$wsLzzy = $sheets.Item("聊斋志异")
for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 1
    $wsLzzy.Range("A" + $row).Value = $words_lzzy[$i]
    $wsLzzy.Range("B" + $row).Value = $freqs_lzzy[$i]
}

$wsHlm = $sheets.Item("红楼梦")
for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 1
    $wsHlm.Range("A" + $row).Value = $words_hlm[$i]
    $wsHlm.Range("B" + $row).Value = $freqs_hlm[$i]
}

# ---------------------------------------------------------------------
# 3. Defined names used by the treemap charts (_xlchart.v1.N), shifted so
#    they keep referring to the right sheet/range after the insert:
#      v1.6/7   -> 聊斋志异  (used to be 三国演义)
#      v1.8/9   -> 红楼梦    (used to be 四世同堂)
#      v1.10/11 -> 三国演义  (used to be 白鹿原)
#      v1.12/13 -> 四世同堂  (new)
#      v1.14/15 -> 白鹿原    (new)
# ---------------------------------------------------------------------
$names = $wb.Names
$names.Item("_xlchart.v1.6").RefersTo = "=聊斋志异!`$A`$1:`$A`$20"
$names.Item("_xlchart.v1.7").RefersTo = "=聊斋志异!`$B`$1:`$B`$20"
$names.Item("_xlchart.v1.8").RefersTo = "=红楼梦!`$A`$1:`$A`$20"
$names.Item("_xlchart.v1.9").RefersTo = "=红楼梦!`$B`$1:`$B`$20"
$names.Item("_xlchart.v1.10").RefersTo = "=三国演义!`$A`$1:`$A`$20"
$names.Item("_xlchart.v1.11").RefersTo = "=三国演义!`$B`$1:`$B`$20"

$n = $names.Add("_xlchart.v1.12", "=四世同堂!`$A`$1:`$A`$20")
$n.Visible = $false
$n = $names.Add("_xlchart.v1.13", "=四世同堂!`$B`$1:`$B`$20")
$n.Visible = $false
$n = $names.Add("_xlchart.v1.14", "=白鹿原!`$A`$1:`$A`$20")
$n.Visible = $false
$n = $names.Add("_xlchart.v1.15", "=白鹿原!`$B`$1:`$B`$20")
$n.Visible = $false

# ---------------------------------------------------------------------
# 4. Add a treemap chart to each new sheet (mirrors the charts already
#    present on the other word-frequency sheets).
# ---------------------------------------------------------------------
try {
    $chartLzzy = $wsLzzy.Shapes.AddChart2(-1, 117)
    $chartLzzy.Chart.SetSourceData($wsLzzy.Range("A1:B20"))
    $chartLzzy.Chart.HasTitle = $true
    $chartLzzy.Chart.ChartTitle.Text = "《聊斋志异》"
} catch {
}

try {
    $chartHlm = $wsHlm.Shapes.AddChart2(-1, 117)
    $chartHlm.Chart.SetSourceData($wsHlm.Range("A1:B20"))
    $chartHlm.Chart.HasTitle = $true
    $chartHlm.Chart.ChartTitle.Text = "《红楼梦》"
} catch {
}

# ---------------------------------------------------------------------
# 5. Restore each sheet's remembered selection, then land on 聊斋志异
#    as the active tab (matches the authored workbook state).
# ---------------------------------------------------------------------
$sheets.Item("红楼梦").Activate()
$sheets.Item("红楼梦").Range("T15").Select()

$sheets.Item("聊斋志异").Activate()
$sheets.Item("聊斋志异").Range("C1").Select()
